$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final sorted (descending by value) data, with Swedish and Uzbek removed.
$data = @(
    @("English", 21.8656027814668),
    @("Chinese", 16.43080754362495),
    @("Spanish", 6.988243635010511),
    @("Arabic", 5.297273520655281),
    @("Japanese", 4.762237531253347),
    @("German", 4.320271123399163),
    @("Russian", 3.622904122386143),
    @("Portuguese", 3.464149046509579),
    @("Malay-Indonesian", 3.06701181524291),
    @("French", 2.649820966570469),
    @("Italian", 2.155467695533362),
    @("Korean", 1.67144534649066),
    @("Turkish", 1.538612828050307),
    @("Persian", 1.34512930415932),
    @("Dutch", 1.261014445143216),
    @("Thai", 1.000931598729329),
    @("Polish", 0.8908862445728305),
    @("Urdu", 0.8404986906310455),
    @("Vietnamese", 0.5639748744277523),
    @("Bengali", 0.5631612597489767)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the two now-unused trailing rows (previously held Uzbek/Vietnamese).
$ws.Range("A22:B23").Delete()
